$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.085.30"
$ws.Range("E2").Value = "  +0.02%  "

$ws.Range("D3").Value = "2.749.64"
$ws.Range("E3").Value = "  -0.07%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "'572.61"
$ws.Range("E5").Value = "  -0.80%  "

$ws.Range("D6").Value = "'159.18"
$ws.Range("E6").Value = "  +0.28%  "

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("E8").Value = "  -1.57%  "

$ws.Range("E9").Value = "  -1.30%  "

$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  +3.84%  "

$ws.Range("E11").Value = "  +1.47%  "

$ws.Range("D12").Value = "'0.385"
$ws.Range("E12").Value = "  -0.87%  "

$ws.Range("D13").Value = "3.237.01"
$ws.Range("E13").Value = "  -0.19%  "

$ws.Range("D14").Value = "'26.69"
$ws.Range("E14").Value = "  -0.28%  "

$ws.Range("D15").Value = "63.722.43"
$ws.Range("E15").Value = "  -0.40%  "

$ws.Range("E16").Value = "  -1.66%  "

$ws.Range("D17").Value = "2.749.81"
$ws.Range("E17").Value = "  -0.31%  "

$ws.Range("D18").Value = "'12.13"
$ws.Range("E18").Value = "  -0.11%  "

$ws.Range("D19").Value = "'4.83"
$ws.Range("E19").Value = "  -1.64%  "

$ws.Range("D20").Value = "'355.40"
$ws.Range("E20").Value = "  -1.34%  "

$ws.Range("D21").Value = "'6.65"
$ws.Range("E21").Value = "  -2.50%  "

$ws.Range("D22").Value = "'0.999"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("D23").Value = "'0.524"
$ws.Range("E23").Value = "  -7.47%  "

$ws.Range("D24").Value = "'64.86"
$ws.Range("E24").Value = "  -2.12%  "

$ws.Range("E25").Value = "  -1.03%  "

$ws.Range("D26").Value = "'8.52"
$ws.Range("E26").Value = "  -0.13%  "

$ws.Range("E27").Value = "  +0.13%  "

$ws.Range("D28").Value = "0.0₃0918"
$ws.Range("E28").Value = "  -1.41%  "

$ws.Range("D29").Value = "'7.32"
$ws.Range("E29").Value = "  +3.66%  "

$ws.Range("E30").Value = "  -0.99%  "

$ws.Range("E31").Value = "  +8.83%  "

$ws.Range("D32").Value = "'167.48"
$ws.Range("E32").Value = "  -0.92%  "

$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").Value = "'4.94"
$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("D34").Value = "'20.18"
$ws.Range("E34").Value = "  -1.41%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.49"
$ws.Range("E35").Value = "  +2.27%  "

$ws.Range("E36").Value = "  +0.01%  "

$ws.Range("D37").Value = "'1.81"
$ws.Range("E37").Value = "  +0.41%  "

$ws.Range("D38").Value = "'0.998"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("D39").Value = "'350.79"
$ws.Range("E39").Value = "  +5.27%  "

$ws.Range("D40").Value = "'6.35"
$ws.Range("E40").Value = "  +4.23%  "

$ws.Range("D41").Value = "'4.15"
$ws.Range("E41").Value = "  -0.88%  "

$ws.Range("D42").Value = "'38.68"
$ws.Range("E42").Value = "  -2.41%  "

$ws.Range("D43").Value = "'22.46"
$ws.Range("E43").Value = "  +2.72%  "

$ws.Range("D44").Value = "'21.51"
$ws.Range("E44").Value = "  -1.57%  "

$ws.Range("D45").Value = "'0.0588"
$ws.Range("E45").Value = "  -0.84%  "

$ws.Range("D46").Value = "'136.51"
$ws.Range("E46").Value = "  +0.32%  "

$ws.Range("D47").Value = "'0.629"
$ws.Range("E47").Value = "  -1.55%  "

$ws.Range("E48").Value = "  -1.15%  "

$ws.Range("E49").Value = "  -2.23%  "

$ws.Range("E50").Value = "  -0.06%  "

$ws.Range("B51").Value = "Maker"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "2.137.80"
$ws.Range("E51").Value = "  +0.98%  "

